$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-formatted cells (column D price strings, etc.) to keep their
# original text representation instead of being auto-coerced to numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.565.93'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.623.40'
$ws.Range('E3').Value = '  -1.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.57'
$ws.Range('E5').Value = '  -0.75%  '
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.25'
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  +1.75%  '
$ws.Range('E10').Value = '  -0.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0891'
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.853.77'
$ws.Range('E12').Value = '  -1.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.612.67'
$ws.Range('E13').Value = '  -1.99%  '
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('E15').Value = '  -2.17%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.34'
$ws.Range('E16').Value = '  +0.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.530.35'
$ws.Range('E17').Value = '  -0.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '231.47'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('E19').Value = '  -0.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.54'
$ws.Range('E20').Value = '  -1.00%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.43'
$ws.Range('E22').Value = '  +2.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.33'
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('E24').Value = '  +6.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.71'
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('E26').Value = '  -0.84%  '
$ws.Range('E27').Value = '  -0.66%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.55'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'BinanceUSD'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('E31').Value = '  -0.52%  '
$ws.Range('E32').Value = '  -0.79%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.464.47'
$ws.Range('E33').Value = '  +1.68%  '
$ws.Range('E34').Value = '  -2.38%  '
$ws.Range('E35').Value = '  -2.75%  '
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.944'
$ws.Range('E37').Value = '  +6.88%  '
$ws.Range('E38').Value = '  +0.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.874'
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('E40').Value = '  -2.75%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('E42').Value = '  -2.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '67.58'
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('E44').Value = '  -2.09%  '
$ws.Range('E45').Value = '  -2.28%  '
$ws.Range('E46').Value = '  -5.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.76'
$ws.Range('E47').Value = '  +1.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.763.98'
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '87.37'
$ws.Range('E49').Value = '  +2.15%  '
$ws.Range('E50').Value = '  -2.17%  '
$ws.Range('E51').Value = '  +1.63%  '
